$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.222.11'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.597.38'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.41'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0605'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.93'
$ws.Range('E10').Value = '  -1.50%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').Value = '1.822.86'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '1.603.49'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.64'
$ws.Range('D17').Value = '26.222.61'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '226.85'
$ws.Range('E18').Value = '  +5.79%  '
$ws.Range('D19').Value = '0.0₃0719'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('E20').Value = '  +3.73%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.23'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.62'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('E28').Value = '  +1.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.35'
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.19'
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').Value = '1.442.22'
$ws.Range('E33').Value = '  +3.86%  '
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.46'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.565'
$ws.Range('E37').Value = '  -3.28%  '
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.819'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.75'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E43').Value = '  -1.02%  '
$ws.Range('D44').Value = '1.735.49'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('E45').Value = '  -1.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.42'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.62'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0500'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('E51').Value = '  -3.08%  '
